# Fruta / hortaliza, semanal
# Insert a new weekly record at row 15 (pushing the existing rows 15-67 down
# to 16-68) for "Femacal de La Calera - Papaya", and populate the new row
# with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 15..67 down to 16..68, creating a blank row 15.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the new weekly entry.
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 45071
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108004
$ws.Range("J15").Value = "Papaya"
$ws.Range("K15").Value = "Cultivar IV Región"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 56
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = "`$/bandeja 10 kilos"
$ws.Range("R15").Value = "Provincia del Elquí"
$ws.Range("S15").Value = 2000
$ws.Range("T15").Value = 10
